$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 707.2222
$ws.Range("I2").Value = 44.166668
$ws.Range("K2").Value = 44.166668
$ws.Range("M2").Value = 68.833332
$ws.Range("H86").Value = 6599
$ws.Range("I86").Value = 7599.5
$ws.Range("K86").Value = 7599.5
$ws.Range("M86").Value = -6476.5
$ws.Range("H89").Value = 6599
$ws.Range("I89").Value = 7599.5
$ws.Range("K89").Value = 37997.5
$ws.Range("M89").Value = -32381.5
$ws.Range("H116").Value = 17625.096
$ws.Range("I116").Value = 4413.778
$ws.Range("K116").Value = 4413.778
$ws.Range("M116").Value = -971.7780000000002
$ws.Range("H137").Value = 23814634
$ws.Range("J137").Value = 6931.407
$ws.Range("L137").Value = 20794.221
$ws.Range("N137").Value = -25894.221
$ws.Range("H138").Value = 4212.037
$ws.Range("I138").Value = 1824.6154
$ws.Range("J138").Value = 6428.9287
$ws.Range("K138").Value = 5473.8462
$ws.Range("L138").Value = 19286.7861
$ws.Range("M138").Value = -333.8462
$ws.Range("N138").Value = -29566.7861

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 231585.06
$ws.Range("I32").Value = 508545.66
$ws.Range("K32").Value = 508545.66
$ws.Range("M32").Value = -508258.66
$ws.Range("H45").Value = 1991
$ws.Range("I45").Value = 1986.5
$ws.Range("K45").Value = 1986.5
$ws.Range("M45").Value = -1609.5
$ws.Range("H61").Value = 1855864.6
$ws.Range("I61").Value = 4074.4146
$ws.Range("J61").Value = 7696126
$ws.Range("K61").Value = 4074.4146
$ws.Range("L61").Value = 7696126
$ws.Range("M61").Value = -3862.4146
$ws.Range("N61").Value = -7696550
$ws.Range("H110").Value = 1372.1
$ws.Range("I110").Value = 1049.5
$ws.Range("J110").Value = 1587.1666
$ws.Range("K110").Value = 1049.5
$ws.Range("L110").Value = 1587.1666
$ws.Range("M110").Value = 995.5
$ws.Range("N110").Value = -5677.1666
$ws.Range("H124").Value = 39400
$ws.Range("J124").Value = 39400
$ws.Range("L124").Value = 39400
$ws.Range("N124").Value = -49220
$ws.Range("H136").Value = 1855864.6
$ws.Range("I136").Value = 4074.4146
$ws.Range("J136").Value = 7696126
$ws.Range("K136").Value = 12223.2438
$ws.Range("L136").Value = 23088378
$ws.Range("M136").Value = -9673.2438
$ws.Range("N136").Value = -23093478

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1699.6111
$ws.Range("I86").Value = 1412.0625
$ws.Range("K86").Value = 1412.0625
$ws.Range("M86").Value = -289.0625
$ws.Range("H89").Value = 1699.6111
$ws.Range("I89").Value = 1412.0625
$ws.Range("K89").Value = 7060.3125
$ws.Range("M89").Value = -1444.3125
$ws.Range("H100").Value = 17999.5
$ws.Range("J100").Value = 17999.5
$ws.Range("L100").Value = 17999.5
$ws.Range("N100").Value = -20163.5
$ws.Range("H105").Value = 6123.273
$ws.Range("I105").Value = 4890.1875
$ws.Range("K105").Value = 4890.1875
$ws.Range("M105").Value = -3143.1875
$ws.Range("H107").Value = 9615987
$ws.Range("I107").Value = 11905359
$ws.Range("J107").Value = 624.8
$ws.Range("K107").Value = 11905359
$ws.Range("L107").Value = 624.8
$ws.Range("M107").Value = -11903439
$ws.Range("N107").Value = -4464.8
$ws.Range("H123").Value = 85000
$ws.Range("J123").Value = 85000
$ws.Range("L123").Value = 85000
$ws.Range("N123").Value = -94800   # LeveProfitHQ cell newly populated (previously absent/blank)
$ws.Range("H134").Value = 2879316.5
$ws.Range("I134").Value = 3895.8235
$ws.Range("K134").Value = 11687.4705
$ws.Range("M134").Value = -9152.470499999999

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 389002.34
$ws.Range("I6").Value = 428574.44
$ws.Range("J6").Value = 250500
$ws.Range("K6").Value = 428574.44
$ws.Range("L6").Value = 250500
$ws.Range("M6").Value = -428461.44
$ws.Range("N6").Value = -250726
$ws.Range("H107").Value = 495.5909
$ws.Range("I107").Value = 497.2
$ws.Range("K107").Value = 497.2
$ws.Range("M107").Value = 1422.8
$ws.Range("H132").Value = 1816.9395
$ws.Range("I132").Value = 1891.826
$ws.Range("K132").Value = 5675.478
$ws.Range("M132").Value = -3145.478

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8451.923000000001
$ws.Range("J68").Value = 8451.923000000001
$ws.Range("L68").Value = 25355.769
$ws.Range("N68").Value = -26977.769
$ws.Range("H71").Value = 8451.923000000001
$ws.Range("J71").Value = 8451.923000000001
$ws.Range("L71").Value = 76067.307
$ws.Range("N71").Value = -84179.307
$ws.Range("H131").Value = 4633705
$ws.Range("I131").Value = 1429.375
$ws.Range("J131").Value = 6949843
$ws.Range("K131").Value = 4288.125
$ws.Range("L131").Value = 20849529
$ws.Range("M131").Value = 751.875
$ws.Range("N131").Value = -20859609
$ws.Range("H137").Value = 6392.394
$ws.Range("I137").Value = 2932.7058
$ws.Range("J137").Value = 10068.3125
$ws.Range("K137").Value = 8798.117400000001
$ws.Range("L137").Value = 30204.9375
$ws.Range("M137").Value = -3698.117400000001
$ws.Range("N137").Value = -40404.9375

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17691.846
$ws.Range("I70").Value = 16932.834
$ws.Range("J70").Value = 19399.625
$ws.Range("K70").Value = 16932.834
$ws.Range("L70").Value = 19399.625
$ws.Range("M70").Value = -16662.834
$ws.Range("N70").Value = -19939.625
$ws.Range("H73").Value = 17691.846
$ws.Range("I73").Value = 16932.834
$ws.Range("J73").Value = 19399.625
$ws.Range("K73").Value = 16932.834
$ws.Range("L73").Value = 19399.625
$ws.Range("M73").Value = -15996.834
$ws.Range("N73").Value = -21271.625
$ws.Range("H122").Value = 6546.5
$ws.Range("J122").Value = 4758.3335
$ws.Range("L122").Value = 14275.0005
$ws.Range("N122").Value = -19175.0005
$ws.Range("H126").Value = 6265.684
$ws.Range("I126").Value = 8327.083000000001
$ws.Range("K126").Value = 24981.249
$ws.Range("M126").Value = -22511.249

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 23199.6
$ws.Range("J38").Value = 23199.6
$ws.Range("L38").Value = 23199.6
$ws.Range("N38").Value = -24019.6
$ws.Range("H93").Value = 2318.5
$ws.Range("I93").Value = 1892.0769
$ws.Range("K93").Value = 1892.0769
$ws.Range("M93").Value = -644.0769
$ws.Range("H122").Value = 3913.9
$ws.Range("J122").Value = 4988
$ws.Range("L122").Value = 14964
$ws.Range("N122").Value = -19864
$ws.Range("H132").Value = 4678425.5
$ws.Range("I132").Value = 8348748.5
$ws.Range("K132").Value = 25046245.5
$ws.Range("M132").Value = -25043715.5

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 10000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()   # LeveProfitHQ cell removed entirely (no longer applicable)
$ws.Range("H81").Value = 4111.091
$ws.Range("I81").Value = 1441.5
$ws.Range("K81").Value = 2883
$ws.Range("M81").Value = -1822
$ws.Range("H84").Value = 4111.091
$ws.Range("I84").Value = 1441.5
$ws.Range("K84").Value = 14415
$ws.Range("M84").Value = -9111
$ws.Range("H126").Value = 2667.8572
$ws.Range("I126").Value = 3188
$ws.Range("K126").Value = 9564
$ws.Range("M126").Value = -7094
$ws.Range("H135").Value = 78357
$ws.Range("J135").Value = 78357
$ws.Range("L135").Value = 78357
$ws.Range("N135").Value = -88497
